$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: add a "Changed from:" comment cell above the renamed column (I) ---
# Copy the formatting used by the other "Changed from:" comment cells (D1/E1) onto I1.
$ws.Range("D1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Changed from: target_protein_lincs_id"

# --- Header row 2: rename the column header and match the D2/E2 style ---
$ws.Range("D2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = "target_protein_center_ids"

# Row 2 grows to the same height as row 1 once it carries a wrapped comment-affected header.
$ws.Rows.Item(2).RowHeight = 68.65

# --- Data: target_protein is now a many-to-many relationship; first row lists two ids ---
$ws.Range("I3").Value = "HMSL201294; HMSL201295"

# --- Misc view state: active cell moves to I5 in the frozen bottom-right pane ---
$ws.Range("I5").Select()

# --- Cosmetic workbook window tab ratio tweak ---
$wb.Windows.Item(1).TabRatio = 0.989
